$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-35, per the regenerated save_data
$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 3
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 4
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 0
    33 = 2
    34 = 2
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
